$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated energy carrier price mark-up values (columns F:AU = years 2010-2051)
# for rows 2, 5, 14 and 17 (historic + future scenario simulation update)

$rowValues2 = @(0.024699999999999993, 0.03500000000000001, 0.04135, 0.05785, 0.0763, 0.07185, 0.0749, 0.0785, 0.0775, 0.08765, 0.102, 0.10005000000000001, 0.05510000000000001, 0.035299999999999984, 0.04, 0.045, 0.05, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055, 0.055)
for ($i = 0; $i -lt $rowValues2.Length; $i++) {
    $ws.Cells.Item(2, 6 + $i).Value = $rowValues2[$i]
}

$rowValues5 = @(0.0049499999999999995, 0.0049499999999999995, 0.0040500000000000015, 0.0040500000000000015, 0.0040500000000000015, 0.003999999999999997, 0.004049999999999998, 0.004049999999999998, 0.0040999999999999995, 0.0050499999999999955, 0.006000000000000002, 0.010199999999999997, 0.011200000000000002, 0.01575, 0.015, 0.014, 0.013, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002, 0.011200000000000002)
for ($i = 0; $i -lt $rowValues5.Length; $i++) {
    $ws.Cells.Item(5, 6 + $i).Value = $rowValues5[$i]
}

$rowValues14 = @(0.0647, 0.0725, 0.07784999999999999, 0.0963, 0.10645000000000002, 0.1049, 0.11134999999999999, 0.1175, 0.1137, 0.11095000000000001, 0.11370000000000001, 0.11220000000000001, 0.0673, 0.04954999999999998, 0.05555, 0.06, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675, 0.0675)
for ($i = 0; $i -lt $rowValues14.Length; $i++) {
    $ws.Cells.Item(14, 6 + $i).Value = $rowValues14[$i]
}

$rowValues17 = @(0.005999999999999998, 0.005899999999999999, 0.005899999999999999, 0.00595, 0.006000000000000002, 0.00595, 0.00595, 0.00595, 0.00595, 0.006100000000000001, 0.0063, 0.010800000000000004, 0.011799999999999998, 0.016, 0.015, 0.014, 0.013, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012, 0.012)
for ($i = 0; $i -lt $rowValues17.Length; $i++) {
    $ws.Cells.Item(17, 6 + $i).Value = $rowValues17[$i]
}

# Update the saved selection / view state on the Sheet1 sheet view
[void]$ws.Range("I35").Select()
